# Add a new "Consume" item "MINYAK MIE SHALLOT OIL" above row 8,
# pushing all subsequent rows down by one (matches the sharedStrings.xml
# diff which appends the new string at the end of the shared-string table,
# and the sheet1.xml diff which shows every row from the old row 8 onward
# shifted down by exactly one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "Consume"
$ws.Range("B8").Value = "MINYAK MIE SHALLOT OIL"

# Leave the same cell selected as in the saved workbook.
$ws.Range("C10").Select() | Out-Null
